$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every row's "distance" value (column B) collapses to the same constant –
# the heatmap side-pixel colorset fix described in the commit message.
$value = 8660.25403784423

for ($i = 1; $i -le 100; $i++) {
    $ws.Cells.Item($i, 2).Value = $value
}
